# "Statistics and to do"
# Adds a small descriptive-statistics block (Average, Standard Deviation, Median,
# Mode, Maximum, Minimum, Range) for the Acceleration column plus a LINEST-based
# linear estimate (slope/intercept) block, then reflows the sheet (auto-fit
# columns, zoom, full-sheet selection) and slides Chart 3 out of the way of the
# new table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("freefall_times_positions")
$ws.Activate()

# ---------------------------------------------------------------------------
# Descriptive statistics on the Acceleration(m/s^2) column (E4:E11)
# ---------------------------------------------------------------------------
$ws.Range("E13").Formula = "=AVERAGE(E4:E11)"
$ws.Range("F13").Value = "Average"

$ws.Range("E14").Formula = "=STDEV.S(E4:E11)"
$ws.Range("F14").Value = "Standard Deviation"

$ws.Range("E16").Formula = "=MEDIAN(E4:E11)"
$ws.Range("F16").Value = "Median"

$ws.Range("E17").Formula = "=MODE.SNGL(E4:E11)"
$ws.Range("F17").Value = "Mode"

$ws.Range("E19").Formula = "=MAX(E4:E11)"
$ws.Range("F19").Value = "Maximum"

$ws.Range("E20").Formula = "=MIN(E4:E11)"
$ws.Range("F20").Value = "Minimum"

$ws.Range("E21").Formula = "=E19-E20"
$ws.Range("F21").Value = "Range"

# ---------------------------------------------------------------------------
# Linear estimate (LINEST) of acceleration vs time
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = "Linear Estimate"
$ws.Range("E22").Value = "Intercept (initial velocity)"
$ws.Range("D22").Value = "Slope (Acceleration)"

$ws.Range("D23:E27").FormulaArray = "=LINEST(D3:D11, B3:B11, TRUE, TRUE)"

# ---------------------------------------------------------------------------
# Reflow: select everything (as after an Auto-fit pass) and resize the
# columns that now hold the wider labels / numbers
# ---------------------------------------------------------------------------
$ws.Cells.Select()

$ws.Columns.Item(1).ColumnWidth = 10.736979166666666
$ws.Columns.Item(2).ColumnWidth = 9.736979166666666
$ws.Columns.Item(3).ColumnWidth = 13.736979166666666
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws.Columns.Item(5).ColumnWidth = 22.307291666666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668

$excel.ActiveWindow.Zoom = 90

# ---------------------------------------------------------------------------
# Move "Chart 3" (the Acceleration chart) to the right so it no longer
# overlaps the new statistics block in columns C:F
# ---------------------------------------------------------------------------
$chart3 = $ws.ChartObjects(3)
$chart3.Left = 892.1665748031496
$chart3.Top = 40.87496062992126
$chart3.Width = 433.39588582677175
$chart3.Height = 216

$ws.Range("A1").Select()
